# Weekly update: a new Pomelo price-report row for
# "Mercado Mayorista Lo Valledor de Santiago" is inserted at row 23,
# pushing the existing rows 23:60 down to 24:61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 23 (shifts rows 23-60 down to 24-61,
# and extends the sheet's used range / dimension to A1:T61).
$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with this week's record.
$ws.Cells.Item(23, 1).Value  = 6
$ws.Cells.Item(23, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value  = "Metropolitana"
$ws.Cells.Item(23, 4).Value  = 45162
$ws.Cells.Item(23, 5).Value  = 13
$ws.Cells.Item(23, 6).Value  = "Fruta"
$ws.Cells.Item(23, 7).Value  = 100102
$ws.Cells.Item(23, 8).Value  = "Cítricos"
$ws.Cells.Item(23, 9).Value  = 100102006
$ws.Cells.Item(23, 10).Value = "Pomelo"
$ws.Cells.Item(23, 11).Value = "Start Ruby"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 15
$ws.Cells.Item(23, 14).Value = 180000
$ws.Cells.Item(23, 15).Value = 180000
$ws.Cells.Item(23, 16).Value = 180000
$ws.Cells.Item(23, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(23, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 19).Value = 514
$ws.Cells.Item(23, 20).Value = 350
